$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("protocoltestcasedetails")

# Update existing Y/N flags that flipped
$ws.Range("D19").Value = "N"
$ws.Range("D20").Value = "N"
$ws.Range("D22").Value = "N"

# Rename testcase28 entry and flip its Y/N flag
$ws.Range("B29").Value = "testcase28_manual_sql_notifications"
$ws.Range("C29").Formula = "=_xlfn.CONCAT(""/app/test/testcases/"",B29,"".xlsx"")"
$ws.Range("D29").Value = "Y"

# Add new row 30 - testcase29
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "testcase29_manual_sql_fullname"
$ws.Range("C30").Formula = "=_xlfn.CONCAT(""/app/test/testcases/"",B30,"".xlsx"")"
$ws.Range("D30").Value = "Y"

# Extend the data validation range to include the new row
$validationRange = $ws.Range("D2:D30")
$validationRange.Validation.Delete()
$validationRange.Validation.Add(3, 1, 1, """Y,N""")
$validationRange.Validation.IgnoreBlank = $true
$validationRange.Validation.InCellDropdown = $true
$validationRange.Validation.ShowInput = $true
$validationRange.Validation.ShowError = $true

# Update the view state to match the saved selection
$ws.Activate()
$ws.Range("E30").Select()
$excel.ActiveWindow.ScrollRow = 11
